# Add a new week's worth of measurements (row 8, date 2025-04-08) below the
# existing data, matching the border/format used by the "interior" data rows,
# and relabel the former last row (row 7) so it matches that same interior
# style now that it is no longer the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 7 was previously styled as the "last row" of the table; now
#        that a row 8 is being appended, row 7 should pick up the regular
#        interior-row border formatting (same as row 2's cells). ---
$null = $ws.Range("B2:M2").Copy()
$null = $ws.Range("B7:M7").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Build the new row 8 by cloning row 6's formatting (date style on
#        column A, regular bordered style on B:M), then filling in values. ---
$null = $ws.Range("A6:M6").Copy()
$null = $ws.Range("A8:M8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(8, 1).Value = 45755
$ws.Cells.Item(8, 2).Value = 36.4
$ws.Cells.Item(8, 3).Value = 42
$ws.Cells.Item(8, 4).Value = 38.9
$ws.Cells.Item(8, 5).Value = 39.200000000000003
$ws.Cells.Item(8, 6).Value = 39.1
$ws.Cells.Item(8, 7).Value = 33.700000000000003
$ws.Cells.Item(8, 8).Value = 32.6
$ws.Cells.Item(8, 9).Value = 39.200000000000003
$ws.Cells.Item(8, 10).Value = 32.6
$ws.Cells.Item(8, 11).Value = 34.5
$ws.Cells.Item(8, 12).Value = 27.8
$ws.Cells.Item(8, 13).Value = 34.799999999999997

# --- 3. Move the active selection to D12 (matches the saved view state). ---
$null = $ws.Range("D12").Select()
